$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 52: was "Metrics running on Africa" summary row, becomes the
#     first row of a new "0_voice" / "0_background" pair (epsilon added
#     to squash -> new metrics rows) ---
$ws.Range("A52").Value = 43377
$ws.Range("B52").Value = 85
$ws.Range("C52").Value = "0_voice"
$ws.Range("D52").Value = [double]"8.8769710273481905E-3"
$ws.Range("E52").Value = -4.4486933776710904
$ws.Range("F52").Value = -2.6658631119058098
$ws.Range("G52").Value = 8.5482997649644705
$ws.Range("H52").Value = [double]"-1.5016279769943E-2"
$ws.Range("I52").Value = "Unet"

# --- Row 53: second half of the merged A52:A53 / B52:B53 pair ---
$ws.Range("C53").Value = "0_background"
$ws.Range("D53").Value = [double]"8.8769710273481905E-3"
$ws.Range("E53").Value = 1.15316280088379
$ws.Range("F53").Value = 5.3171892248921999
$ws.Range("G53").Value = 8.5446161713002304
$ws.Range("H53").Value = [double]"2.0442723267735299E-2"

# --- Row 54: was "Final test (after 9 epochs)..." summary row, becomes
#     new "0_voice" metrics row ---
$ws.Range("A54").Value = 43378
$ws.Range("B54").Value = 89
$ws.Range("C54").Value = "0_voice"
$ws.Range("D54").Value = [double]"6.0007950481425097E-3"
$ws.Range("E54").Value = 0.53187985334754795
$ws.Range("F54").Value = 6.2616599956293602
$ws.Range("G54").Value = 5.3511069949812704
$ws.Range("H54").Value = 4.9655569512486801

# --- Row 55: second half of the merged A54:A55 / B54:B55 pair ---
$ws.Range("C55").Value = "0_background"
$ws.Range("D55").Value = [double]"6.0007950481425097E-3"
$ws.Range("E55").Value = 4.78322596010337
$ws.Range("F55").Value = 13.387147861798701
$ws.Range("G55").Value = 6.5279267189934798
$ws.Range("H55").Value = 3.6505058824873098

# --- Row 56: brand new standalone row carrying what used to be row 52's
#     description text, now reworded to "Running on Asia" ---
$ws.Range("A56").Value = 43378
$ws.Range("A56").NumberFormat = "d-mmm"
$ws.Range("B56").Value = 88
$ws.Range("C56").Value = "Running on Asia"

# --- View state: scroll position / selection moved ---
[void]$ws.Range("I58").Select()
